$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 3-7, column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 117
$ws1.Range("F4").Value = 151
$ws1.Range("F5").Value = 3067
$ws1.Range("F6").Value = 311
$ws1.Range("F7").Value = 408

# Sheet "全部类型" (All types) - rows 3,4,5,6,9, column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 117
$ws4.Range("F4").Value = 151
$ws4.Range("F5").Value = 3067
$ws4.Range("F6").Value = 311
$ws4.Range("F9").Value = 408
